$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''68.622.57'
$ws.Range("E2").Value = '  -1.37%  '

$ws.Range("D3").Value = '''2.450.93'
$ws.Range("E3").Value = '  -1.78%  '

$ws.Range("E4").Value = '  -0.01%  '

$ws.Range("D5").Value = '''556.79'
$ws.Range("E5").Value = '  -2.38%  '

$ws.Range("E6").Value = '  -3.35%  '

$ws.Range("E7").Value = '  +0.02%  '

$ws.Range("D8").Value = '''0.502'
$ws.Range("E8").Value = '  -1.61%  '

$ws.Range("D9").Value = '''2.448.86'
$ws.Range("E9").Value = '  -1.84%  '

$ws.Range("E10").Value = '  -7.15%  '

$ws.Range("E11").Value = '  -1.02%  '

$ws.Range("D12").Value = '''0.332'
$ws.Range("E12").Value = '  -5.54%  '

$ws.Range("D13").Value = '''4.79'
$ws.Range("E13").Value = '  -1.84%  '

$ws.Range("D14").Value = '''2.897.67'
$ws.Range("E14").Value = '  -1.83%  '

$ws.Range("D15").Value = '''68.511.21'
$ws.Range("E15").Value = '  -1.20%  '

$ws.Range("E16").Value = '  -4.35%  '

$ws.Range("D17").Value = '''23.32'
$ws.Range("E17").Value = '  -3.85%  '

$ws.Range("D18").Value = '''2.417.87'
$ws.Range("E18").Value = '  -3.96%  '

$ws.Range("D19").Value = '''10.68'
$ws.Range("E19").Value = '  -5.01%  '

$ws.Range("D20").Value = '''340.43'
$ws.Range("E20").Value = '  -4.56%  '

$ws.Range("E21").Value = '  -5.84%  '

$ws.Range("E22").Value = '  -3.08%  '

$ws.Range("D23").Value = '''6.03'
$ws.Range("E23").Value = '  -0.78%  '

$ws.Range("E24").Value = '  +0.01%  '

$ws.Range("E25").Value = '  -2.72%  '

$ws.Range("D26").Value = '''66.56'
$ws.Range("E26").Value = '  -4.13%  '

$ws.Range("E27").Value = '  -5.31%  '

$ws.Range("D28").Value = '''2.574.59'
$ws.Range("E28").Value = '  -1.79%  '

$ws.Range("D29").Value = '''0.999'
$ws.Range("E29").Value = '  +0.16%  '

$ws.Range("E30").Value = '  -6.31%  '

$ws.Range("D31").Value = '''0.0₃0811'
$ws.Range("E31").Value = '  -7.20%  '

$ws.Range("D32").Value = '''7.12'
$ws.Range("E32").Value = '  -6.43%  '

$ws.Range("B33").Value = 'FirstDigitalUSD'
$ws.Range("C33").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D33").Value = '''0.999'
$ws.Range("E33").Value = '  -0.04%  '

$ws.Range("B34").Value = 'Bittensor'
$ws.Range("C34").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D34").Value = '''432.68'
$ws.Range("E34").Value = '  -1.16%  '

$ws.Range("E35").Value = '  -5.28%  '

$ws.Range("E36").Value = '  -6.18%  '

$ws.Range("D37").Value = '''157.39'
$ws.Range("E37").Value = '  +1.58%  '

$ws.Range("D38").Value = '''19.03'
$ws.Range("E38").Value = '  -0.20%  '

$ws.Range("E39").Value = '  -0.02%  '

$ws.Range("E40").Value = '  -3.35%  '

$ws.Range("D41").Value = '''17.78'
$ws.Range("E41").Value = '  -2.21%  '

$ws.Range("E42").Value = '  -4.09%  '

$ws.Range("D43").Value = '''4.38'
$ws.Range("E43").Value = '  -4.55%  '

$ws.Range("D44").Value = '''37.41'
$ws.Range("E44").Value = '  -0.90%  '

$ws.Range("E45").Value = '  -7.56%  '

$ws.Range("E46").Value = '  +2.48%  '

$ws.Range("E47").Value = '  -6.11%  '

$ws.Range("D48").Value = '''131.73'
$ws.Range("E48").Value = '  -4.90%  '

$ws.Range("D49").Value = '''3.33'
$ws.Range("E49").Value = '  -2.95%  '

$ws.Range("D50").Value = '''0.0712'
$ws.Range("E50").Value = '  -1.46%  '

$ws.Range("D51").Value = '''0.480'
$ws.Range("E51").Value = '  -5.17%  '
